# Add data for 2022-07-07
# Updates citywide totals, by-neighborhood totals, and per-neighborhood
# breakdowns in the violent-crime-full-year workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 3515
$ws.Range("H3").Value = 8347
$ws.Range("I3").Value = 3662
$ws.Range("H4").Value = 1667
$ws.Range("I4").Value = 856
$ws.Range("I5").Value = 341
$ws.Range("I6").Value = 4100
$ws.Range("H7").Value = 25977
$ws.Range("I7").Value = 12474

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I7").Value = 408
$ws.Range("I8").Value = 760
$ws.Range("I19").Value = 330
$ws.Range("I20").Value = 307
$ws.Range("I22").Value = 35
$ws.Range("I23").Value = 117
$ws.Range("I26").Value = 16
$ws.Range("I27").Value = 111
$ws.Range("I29").Value = 812
$ws.Range("I33").Value = 561
$ws.Range("I37").Value = 394
$ws.Range("I41").Value = 56
$ws.Range("I42").Value = 431
$ws.Range("I43").Value = 110
$ws.Range("I47").Value = 85
$ws.Range("I50").Value = 57
$ws.Range("H52").Value = 525
$ws.Range("I52").Value = 274
$ws.Range("I54").Value = 280
$ws.Range("H63").Value = 207
$ws.Range("I63").Value = 50
$ws.Range("I65").Value = 277
$ws.Range("I67").Value = 483
$ws.Range("I73").Value = 108
$ws.Range("I76").Value = 191
$ws.Range("I77").Value = 69
$ws.Range("I78").Value = 177
$ws.Range("I79").Value = 325
$ws.Range("I83").Value = 250
$ws.Range("I84").Value = 109
$ws.Range("I85").Value = 577
$ws.Range("I88").Value = 114
$ws.Range("I89").Value = 141
$ws.Range("I90").Value = 157
$ws.Range("I91").Value = 152
$ws.Range("I93").Value = 68
$ws.Range("I98").Value = 80
$ws.Range("I99").Value = 232
$ws.Range("H101").Value = 25977
$ws.Range("I101").Value = 12474

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I3").Value = 233
$ws.Range("I6").Value = 144
$ws.Range("I7").Value = 577

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I2").Value = 78
$ws.Range("H4").Value = 29
$ws.Range("I6").Value = 65
$ws.Range("H7").Value = 525
$ws.Range("I7").Value = 274

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I3").Value = 212
$ws.Range("I7").Value = 760

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 139
$ws.Range("I4").Value = 22
$ws.Range("I6").Value = 105
$ws.Range("I7").Value = 408

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I5").Value = 6
$ws.Range("I7").Value = 141

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 126
$ws.Range("I6").Value = 112
$ws.Range("I7").Value = 394

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I3").Value = 85
$ws.Range("I7").Value = 232

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 116
$ws.Range("I3").Value = 171
$ws.Range("I6").Value = 161
$ws.Range("I7").Value = 483

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 109

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I3").Value = 76
$ws.Range("I5").Value = 13
$ws.Range("I6").Value = 86
$ws.Range("I7").Value = 277

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I3").Value = 96
$ws.Range("I7").Value = 250

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 202
$ws.Range("I7").Value = 561

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 54
$ws.Range("I4").Value = 18
$ws.Range("I6").Value = 143
$ws.Range("I7").Value = 280

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 244
$ws.Range("I3").Value = 281
$ws.Range("I7").Value = 812

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 92
$ws.Range("I7").Value = 330

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I4").Value = 21
$ws.Range("I6").Value = 81
$ws.Range("I7").Value = 191

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("I3").Value = 19
$ws.Range("I7").Value = 56

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 117
$ws.Range("I3").Value = 148
$ws.Range("I6").Value = 113
$ws.Range("I7").Value = 431

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I2").Value = 36
$ws.Range("I7").Value = 177

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I3").Value = 40
$ws.Range("I7").Value = 117

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I3").Value = 53
$ws.Range("I7").Value = 152

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 94
$ws.Range("I3").Value = 106
$ws.Range("I6").Value = 97
$ws.Range("I7").Value = 325

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 88
$ws.Range("I7").Value = 307

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("I6").Value = 26
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I6").Value = 32
$ws.Range("I7").Value = 85

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("I5").Value = 2
$ws.Range("I7").Value = 80

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 57

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("I6").Value = 10
$ws.Range("I7").Value = 16

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I2").Value = 39
$ws.Range("I7").Value = 108

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I2").Value = 28
$ws.Range("I7").Value = 114

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I2").Value = 28
$ws.Range("I7").Value = 111

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I6").Value = 54
$ws.Range("I7").Value = 157

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I6").Value = 65
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I4").Value = 2
$ws.Range("I7").Value = 35

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 69
